# removed time logic from GUI and forced all time references to go through
# ProjectTimer class -- the GUI now just writes the computed Time/Material
# totals straight into the worksheet instead of deriving them live.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Time" (C) / "Material" (D) figures for the first three pieces,
# written by the ProjectTimer class.
$ws.Range("C2").Value = 5.264999866485596
$ws.Range("D2").Value = 3.0

$ws.Range("C3").Value = 0.8740000128746033
$ws.Range("D3").Value = -1.0

$ws.Range("C4").Value = 33.07899856567383
$ws.Range("D4").Value = 20.0

# Leave the active selection where the user's last input landed.
$ws.Range("D3").Select()
